# Rename worksheet "Property1" -> "DataNode"
# (unifying the DataNode / DataTable / Entity naming convention per the commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Preserve the author's final cell selection on the sheet (W37), as recorded
# in the saved view state of the edited workbook.
[void]$ws.Range("W37").Select()
